$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: new column I
$ws.Range('I1').Value = 'Other found locations'

# Row 2: update Authors (E2), add Other found locations (I2)
$ws.Range('E2').Value = '[Hamideh%Abbaspour Kasgari%NULL%1,  Siavash%Moradi%NULL%1,  Amir Mohammad%Shabani%NULL%1,  Farhang%Babamahmoodi%NULL%2,  Farhang%Babamahmoodi%NULL%0,  Ali Reza%Davoudi Badabi%NULL%1,  Lotfollah%Davoudi%NULL%1,  Ahmad%Alikhani%NULL%1,  Akbar%Hedayatizadeh Omran%NULL%1,  Majid%Saeedi%NULL%2,  Shahin%Merat%NULL%3,  Hannah%Wentzel%NULL%6,  Hannah%Wentzel%NULL%0,  Anna%Garratt%NULL%2,  Jacob%Levi%NULL%3,  Bryony%Simmons%NULL%3,  Andrew%Hill%NULL%4,  Hafez%Tirgar Fakheri%hafezfakheri@gmail.com%1]'
$ws.Range('I2').Value = '_PMC'

# Row 3: add Other found locations (I3) as an empty string
$ws.Range('I3').Value = "'"
$ws.Range('I3').ClearFormats()

# Row 4: update Authors (E4), add Other found locations (I4)
$ws.Range('E4').Value = '[Gholamali%Eslami%NULL%1,  Sajedeh%Mousaviasl%NULL%1,  Esmat%Radmanesh%NULL%1,  Saeed%Jelvay%NULL%1,  Saeid%Bitaraf%NULL%1,  Bryony%Simmons%NULL%0,  Hannah%Wentzel%NULL%0,  Andrew%Hill%NULL%0,  Anahita%Sadeghi%NULL%1,  James%Freeman%NULL%1,  Shokrollah%Salmanzadeh%NULL%1,  Hani%Esmaeilian%NULL%1,  Morteza%Mobarak%NULL%1,  Ramin%Tabibi%NULL%1,  Amir Hosein%Jafari Kashi%NULL%1,  Zahra%Lotfi%NULL%1,  Seyed Mehdi%Talebzadeh%NULL%1,  Aseni%Wickramatillake%NULL%1,  Mahboobeh%Momtazan%NULL%1,  Majid%Hajizadeh Farsani%NULL%1,  Sedigheh%Marjani%NULL%1,  Sara%Mobarak%s.mobarak@abadanums.ac.ir%1]'
$ws.Range('I4').Value = '_PMC'

# Row 5: update Authors (E5), add Other found locations (I5)
$ws.Range('E5').Value = '[Fatemeh%Roozbeh%NULL%1,  Majid%Saeedi%NULL%0,  Reza%Alizadeh-Navaei%NULL%1,  Akbar%Hedayatizadeh-Omran%NULL%1,  Shahin%Merat%NULL%0,  Hannah%Wentzel%NULL%0,  Hannah%Wentzel%NULL%0,  Jacob%Levi%NULL%0,  Andrew%Hill%NULL%0,  Amir%Shamshirian%shamshirian.amir@gmail.com%1]'
$ws.Range('I5').Value = '_PMC'

# Row 6: update Authors (E6), add Other found locations (I6)
$ws.Range('E6').Value = '[Anahid%Nourian%NULL%1,  Hossein%Khalili%NULL%1,  Zahra%Ahmadinejad%NULL%1,  Hamid%Emadi Kouchak%NULL%1,  Sirous%Jafari%NULL%1,  Sayed Ali%Dehghan Manshadi%NULL%1,  Mehrnaz%Rasolinejad%NULL%1,  Abbas%Kebriaeezadeh%NULL%1]'
$ws.Range('I6').Value = '_PMC'

# Row 7: update Authors (E7), add Other found locations (I7)
$ws.Range('E7').Value = '[Babak%Sayad%NULL%1,  Reza%Khodarahmi%rkhodarahmi@mbrc.ac.ir%2,  Reza%Khodarahmi%rkhodarahmi@mbrc.ac.ir%0,  Farid%Najafi%NULL%1,  Ronak%Miladi%NULL%1,  Zeinab%Mohseni Afshar%NULL%1,  Feizollah%Mansouri%NULL%1,  Zohreh%Rahimi%NULL%1,  Maria%Shirvani%NULL%1,  Mehdi%Salimi%NULL%1,  Siavash%Vaziri%NULL%1,  Alireza%Janbakhsh%NULL%1,  Fatemeh%Khosravi Shadmani%NULL%1,  Arezoo%Bozorgomid%NULL%1,  Mohammad Hossein%Zamanian%NULL%1,  Mandana%Afsharian%NULL%1]'
$ws.Range('I7').Value = '_PMC'

# Row 8: update Authors (E8), add Other found locations (I8)
$ws.Range('E8').Value = '[Anahita%Sadeghi%NULL%1,  Ali%Ali Asgari%NULL%2,  Ali%Ali Asgari%NULL%0,  Alireza%Norouzi%NULL%1,  Zahedin%Kheiri%NULL%1,  Amir%Anushirvani%NULL%1,  Mahnaz%Montazeri%NULL%1,  Hadiseh%Hosamirudsai%NULL%1,  Shirin%Afhami%NULL%1,  Elham%Akbarpour%NULL%1,  Rasoul%Aliannejad%NULL%2,  Rasoul%Aliannejad%NULL%0,  Amir Reza%Radmard%NULL%1,  Amir H%Davarpanah%NULL%1,  Jacob%Levi%NULL%0,  Hannah%Wentzel%NULL%0,  Ambar%Qavi%NULL%1,  Anna%Garratt%NULL%0,  Bryony%Simmons%NULL%0,  Andrew%Hill%NULL%0,  Shahin%Merat%merat@tums.ac.ir%0]'
$ws.Range('I8').Value = '_PMC'

# Row 9: update Authors (E9), add Other found locations (I9)
$ws.Range('E9').Value = '[Sherif%Abbass%sherif533@hotmail.com%1,  Ehab%Kamal%NULL%2,  Ehab%Kamal%NULL%0,  Mohsen%Salama%NULL%2,  Mohsen%Salama%NULL%0,  Tary%Salman%NULL%1,  Alyaa%Sabry%NULL%1,  Wael%Abdel‐Razek%NULL%1,  Sherine%Helmy%NULL%1,  Ahmed%Abdelgwad%NULL%1,  Neamt%Sakr%NULL%1,  Mohamed%Elgazzar%NULL%1,  Mohamed%Einar%NULL%1,  Mahmoud%Farouk%NULL%1,  Mounir%Saif%NULL%1,  Ismail%Shehab%NULL%1,  Eman%El‐hosieny%NULL%1,  Mai%Mansour%NULL%1,  Doaa%Mahdi%NULL%1,  El‐Sayed%Tharwa%NULL%1,  Mostafa%Salah%NULL%1,  Ola%Elrouby%NULL%1,  Imam%Waked%NULL%1]'
$ws.Range('I9').Value = '_PMC'

